# "2006 Monthly Time Charts" — insert a new "monthly-time-chart" worksheet
# between "location" and "classification", carrying the 12 monthly larceny
# counts, formatted with a thousands-separator number format.

$wb = $excel.ActiveWorkbook

# Insert the new sheet right after "location" so the final order is:
# location, monthly-time-chart, classification, day-of-week.
$afterSheet = $wb.Worksheets.Item("location")
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $afterSheet)
$ws.Name = "monthly-time-chart"

# Header row: Month, January .. December (row formatted out to column Z,
# matching the rest of the workbook's charts).
$months = @("January","February","March","April","May","June","July","August","September","October","November","December")

$ws.Range("A1:Z1").Font.Name = "Arial"
$ws.Range("A1:Z1").Font.Size = 10

$ws.Range("A1").Value = "Month"
for ($i = 0; $i -lt $months.Length; $i++) {
    $col = [char]([int][char]'B' + $i)
    $ws.Range("$col`1").Value = $months[$i]
}

# Data row: label + monthly larceny counts
$values = @(11929, 10250, 11574, 11838, 13026, 13271, 13093, 13360, 12560, 12283, 11421, 11249)

$ws.Range("A2").Value = "Number of Larcenies"
for ($i = 0; $i -lt $values.Length; $i++) {
    $col = [char]([int][char]'B' + $i)
    $ws.Range("$col`2").Value = $values[$i]
}

# Number format (#,###) for the monthly values
$ws.Range("B2:M2").NumberFormat = "#,###"

# Column width to match the other charts' label column
$ws.Columns.Item(1).ColumnWidth = 34.9

# Leave the workbook focused back on the first sheet, as it was originally.
$wb.Worksheets.Item("location").Activate()
